# Update the player roster table: several players' position/team info
# changed, the roster was re-ordered, and one new player row was appended.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired contents for A2:C19 (header row A1:C1 stays untouched).
$targetData = @(
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Russell Westbrook", "PG,SG", "Denver Nuggets"),
    @("Payton Pritchard", "PG,SG", "Boston Celtics"),
    @("Dejounte Murray", "PG,SG", "New Orleans Pelicans"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Ayo Dosunmu", "PG,SG,SF", "Chicago Bulls"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Jakob Poeltl", "C", "Toronto Raptors"),
    @("Stephon Castle", "PG,SG", "San Antonio Spurs"),
    @("Pascal Siakam", "SF,PF,C", "Indiana Pacers"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Ty Jerome", "PG,SG", "Cleveland Cavaliers"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder"),
    @("Jalen Suggs", "PG,SG", "Orlando Magic"),
    @("Khris Middleton", "SF", "Milwaukee Bucks")
)

for ($i = 0; $i -lt $targetData.Count; $i++) {
    $row = $i + 2
    $vals = $targetData[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
}
